# diary update week 2/6
# Turn the first template placeholder row (row 18) into a real diary entry,
# matching the format/styling already used by rows 12/14/16 (date, time,
# Class, Attend lecture, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats, fonts, fills, wrap text, row height)
# from row 14 -- which already has the exact style pattern this new entry
# needs -- down onto row 18, then overwrite row 18's values/text.
$ws.Range("A14:G14").Copy()
$ws.Range("A18:G18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(18).RowHeight = 72.35

$ws.Range("A18").Value = "02/06/20"
$ws.Range("B18").Value = "5:00:00 PM"

$ws.Range("C18").Value = "Class"
$ws.Range("D18").Value = "Attend lecture"
$ws.Range("E18").Value = "Learned about and practiced mental simulation"
$ws.Range("F18").Value = "Using examples can be a more effective way of simulating than simply reading through the code. You should also be careful to verify that code works the way it appears to."
$ws.Range("G18").Value = "Tricky but neat"

# Move the view the way the author's window ended up: scrolled back toward
# the top of the sheet, with the newly-added row's last cell selected.
$ws.Range("G19").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
